# Scheduled-runner update: refresh cached market-board pricing figures
# (currentAveragePrice / NQ / HQ / Leve price / profit columns H:N) for a
# handful of leve rows across several crafting-sheet tabs. Values below
# are the refreshed figures pulled by the runner; some profit cells that
# are no longer computable for a row are cleared outright.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 91.666664
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H43").Value = 4125
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 4500
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 4500
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -4638
$ws.Range("H116").Value = 4997.6665
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4997.6665
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4997.6665
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -11881.6665
$ws.Range("H132").Value = 4948.1055
$ws.Range("I132").Value = 2216.4285
$ws.Range("J132").Value = 12596.8
$ws.Range("K132").Value = 6649.2855
$ws.Range("L132").Value = 37790.39999999999
$ws.Range("M132").Value = -4119.2855
$ws.Range("N132").Value = -42850.39999999999
$ws.Range("H135").Value = 36066.668
$ws.Range("I135").Value = 36066.668
$ws.Range("K135").Value = 324600.012
$ws.Range("M135").Value = -322065.012
$ws.Range("H137").Value = 8286.726000000001
$ws.Range("I137").Value = 1807.3334
$ws.Range("J137").Value = 17258.191
$ws.Range("K137").Value = 5422.0002
$ws.Range("L137").Value = 51774.573
$ws.Range("M137").Value = -2872.0002
$ws.Range("N137").Value = -56874.573
$ws.Range("H138").Value = 3715.6
$ws.Range("I138").Value = 3327.6428
$ws.Range("J138").Value = 3866.4722
$ws.Range("K138").Value = 9982.928400000001
$ws.Range("L138").Value = 11599.4166
$ws.Range("M138").Value = -4842.928400000001
$ws.Range("N138").Value = -21879.4166
$ws.Range("H141").Value = 3626.4707
$ws.Range("I141").Value = 3303.5715
$ws.Range("J141").Value = 5133.3335
$ws.Range("K141").Value = 9910.7145
$ws.Range("L141").Value = 15400.0005
$ws.Range("M141").Value = -4730.7145
$ws.Range("N141").Value = -25760.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5020.3335
$ws.Range("I32").Value = 2227.3057
$ws.Range("K32").Value = 2227.3057
$ws.Range("M32").Value = -1940.3057
$ws.Range("H47").Value = 24500
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H61").Value = 1480983.2
$ws.Range("I61").Value = 3673.5
$ws.Range("J61").Value = 2465856.2
$ws.Range("K61").Value = 3673.5
$ws.Range("L61").Value = 2465856.2
$ws.Range("M61").Value = -3461.5
$ws.Range("N61").Value = -2466280.2
$ws.Range("H74").Value = 13282.4
$ws.Range("I74").Value = 4308.5713
$ws.Range("K74").Value = 4308.5713
$ws.Range("M74").Value = -3434.5713
$ws.Range("H77").Value = 13282.4
$ws.Range("I77").Value = 4308.5713
$ws.Range("K77").Value = 21542.8565
$ws.Range("M77").Value = -17174.8565
$ws.Range("H136").Value = 1480983.2
$ws.Range("I136").Value = 3673.5
$ws.Range("J136").Value = 2465856.2
$ws.Range("K136").Value = 11020.5
$ws.Range("L136").Value = 7397568.600000001
$ws.Range("M136").Value = -8470.5
$ws.Range("N136").Value = -7402668.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1843.6428
$ws.Range("I105").Value = 1600.8462
$ws.Range("K105").Value = 1600.8462
$ws.Range("M105").Value = 146.1538
$ws.Range("H107").Value = 2301.2222
$ws.Range("I107").Value = 2244.4285
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 2244.4285
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -324.4285
$ws.Range("N107").Value = -6340
$ws.Range("H134").Value = 15992.875
$ws.Range("I134").Value = 10712.741
$ws.Range("K134").Value = 32138.223
$ws.Range("M134").Value = -29603.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 69619.5
$ws.Range("J48").Value = 69619.5
$ws.Range("L48").Value = 69619.5
$ws.Range("N48").Value = -70571.5
$ws.Range("H94").Value = 24503
$ws.Range("I94").Value = 24503
$ws.Range("K94").Value = 24503
$ws.Range("M94").Value = -24052
$ws.Range("H99").Value = 6326.4614
$ws.Range("I99").Value = 5461.7856
$ws.Range("J99").Value = 7335.25
$ws.Range("K99").Value = 5461.7856
$ws.Range("L99").Value = 7335.25
$ws.Range("M99").Value = -3963.7856
$ws.Range("N99").Value = -10331.25
$ws.Range("H126").Value = 6326.4614
$ws.Range("I126").Value = 5461.7856
$ws.Range("J126").Value = 7335.25
$ws.Range("K126").Value = 16385.3568
$ws.Range("L126").Value = 22005.75
$ws.Range("M126").Value = -13915.3568
$ws.Range("N126").Value = -26945.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H56").Value = 125006024
$ws.Range("I56").Value = 125006024
$ws.Range("K56").Value = 125006024
$ws.Range("M56").Value = -125005494
$ws.Range("H131").Value = 1409.6875
$ws.Range("I131").Value = 865.5
$ws.Range("J131").Value = 1487.4286
$ws.Range("K131").Value = 2596.5
$ws.Range("L131").Value = 4462.2858
$ws.Range("M131").Value = 2443.5
$ws.Range("N131").Value = -14542.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 999
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -811
$ws.Range("N46").ClearContents()
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 82105
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 82105
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -83393
$ws.Range("H56").Value = 10949.75
$ws.Range("I56").Value = 10949.75
$ws.Range("K56").Value = 10949.75
$ws.Range("M56").Value = -10258.75
$ws.Range("H99").Value = 28382.572
$ws.Range("I99").Value = 28382.572
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 28382.572
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -25387.572
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 3113.9524
$ws.Range("I100").Value = 2758.25
$ws.Range("J100").Value = 3588.2222
$ws.Range("K100").Value = 2758.25
$ws.Range("L100").Value = 3588.2222
$ws.Range("M100").Value = -2217.25
$ws.Range("N100").Value = -4670.2222
$ws.Range("H132").Value = 1460907.8
$ws.Range("I132").Value = 3809.5715
$ws.Range("J132").Value = 2594206.2
$ws.Range("K132").Value = 11428.7145
$ws.Range("L132").Value = 7782618.600000001
$ws.Range("M132").Value = -8898.7145
$ws.Range("N132").Value = -7787678.600000001
$ws.Range("H136").Value = 974381.0600000001
$ws.Range("J136").Value = 1101459.8
$ws.Range("L136").Value = 3304379.4
$ws.Range("N136").Value = -3309479.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5084.2593
$ws.Range("I122").Value = 3611.3125
$ws.Range("J122").Value = 7226.727
$ws.Range("K122").Value = 10833.9375
$ws.Range("L122").Value = 21680.181
$ws.Range("M122").Value = -8383.9375
$ws.Range("N122").Value = -26580.181
